$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) date values for rows 2-8 from 2023-10-09 (45208) to 2023-10-13 (45212)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45212
}
